# "Generate Report for Handback"
#
# This updates the localization-status report to reflect that the
# handback has completed and is in sync with en-US:
#   - Overview!C2/C3 status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (shared string is reused by both
#     rows, so a single text edit covers both).
#   - zh-cn / de-de detail sheets get their "Latest Target File" (F) and
#     "Latest Handback File" (G) columns populated (with hyperlinks) for
#     both data rows, and the "Latest Handback DateTime" (H) timestamps
#     are stamped with the real handback time.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: status text -----------------------------------
# Columns: A=File Name, B=zh-cn status, C=de-de status, D=Latest Handoff Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet -----------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("F2").Value = "898253a8-a7e4-4d01-b785-db9daf822a06.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e3b12d17401494a259b2b0bbd15128741d23416b/e2e/898253a8-a7e4-4d01-b785-db9daf822a06.md",
    "",
    "",
    "898253a8-a7e4-4d01-b785-db9daf822a06.md"
) | Out-Null

$wsZhCn.Range("G2").Value = "898253a8-a7e4-4d01-b785-db9daf822a06.0dc51ea33dec87553b4fc4e71873925878493b30.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2f953043dad3a1ab1316429d32f1f91fbfab23d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/898253a8-a7e4-4d01-b785-db9daf822a06.0dc51ea33dec87553b4fc4e71873925878493b30.zh-cn.xlf",
    "",
    "",
    "898253a8-a7e4-4d01-b785-db9daf822a06.0dc51ea33dec87553b4fc4e71873925878493b30.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("H2").Value = "2016-03-13 17:05:59"

$wsZhCn.Range("F3").Value = "b4576177-f78d-4cb5-a183-23e52148f132.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e3b12d17401494a259b2b0bbd15128741d23416b/e2e/b4576177-f78d-4cb5-a183-23e52148f132.md",
    "",
    "",
    "b4576177-f78d-4cb5-a183-23e52148f132.md"
) | Out-Null

$wsZhCn.Range("G3").Value = "b4576177-f78d-4cb5-a183-23e52148f132.b30ea201d95720779aa937aaefdefa5f1c6fd288.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2f953043dad3a1ab1316429d32f1f91fbfab23d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b4576177-f78d-4cb5-a183-23e52148f132.b30ea201d95720779aa937aaefdefa5f1c6fd288.zh-cn.xlf",
    "",
    "",
    "b4576177-f78d-4cb5-a183-23e52148f132.b30ea201d95720779aa937aaefdefa5f1c6fd288.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("H3").Value = "2016-03-13 17:05:59"

# ---- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("F2").Value = "898253a8-a7e4-4d01-b785-db9daf822a06.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e3b12d17401494a259b2b0bbd15128741d23416b/e2e/898253a8-a7e4-4d01-b785-db9daf822a06.md",
    "",
    "",
    "898253a8-a7e4-4d01-b785-db9daf822a06.md"
) | Out-Null

$wsDeDe.Range("G2").Value = "898253a8-a7e4-4d01-b785-db9daf822a06.0dc51ea33dec87553b4fc4e71873925878493b30.de-de.xlf"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cef05287ee3f17324125c53109364eca9e17010/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/898253a8-a7e4-4d01-b785-db9daf822a06.0dc51ea33dec87553b4fc4e71873925878493b30.de-de.xlf",
    "",
    "",
    "898253a8-a7e4-4d01-b785-db9daf822a06.0dc51ea33dec87553b4fc4e71873925878493b30.de-de.xlf"
) | Out-Null

$wsDeDe.Range("H2").Value = "2016-03-13 17:06:06"

$wsDeDe.Range("F3").Value = "b4576177-f78d-4cb5-a183-23e52148f132.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e3b12d17401494a259b2b0bbd15128741d23416b/e2e/b4576177-f78d-4cb5-a183-23e52148f132.md",
    "",
    "",
    "b4576177-f78d-4cb5-a183-23e52148f132.md"
) | Out-Null

$wsDeDe.Range("G3").Value = "b4576177-f78d-4cb5-a183-23e52148f132.b30ea201d95720779aa937aaefdefa5f1c6fd288.de-de.xlf"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cef05287ee3f17324125c53109364eca9e17010/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b4576177-f78d-4cb5-a183-23e52148f132.b30ea201d95720779aa937aaefdefa5f1c6fd288.de-de.xlf",
    "",
    "",
    "b4576177-f78d-4cb5-a183-23e52148f132.b30ea201d95720779aa937aaefdefa5f1c6fd288.de-de.xlf"
) | Out-Null

$wsDeDe.Range("H3").Value = "2016-03-13 17:06:06"
